# Add "NA" values under the duplicate_image_filename column (E) for the
# data rows of the stimuli table (rows 2-21).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 5).Value = "NA"
}
